# The commit swaps the presentation's theme colour palette from the
# "Integral" scheme (dk1=000000, lt1=FFFFFF, dk2=455F51, lt2=E3DED1,
# accent1=99CB38 ... folHlink=B26B02) over to the stock "Office" colour
# scheme (dk2=44546A, lt2=E7E6E6, accent1=5B9BD5 ... folHlink=954F72) -
# i.e. the deck's theme (ppt/theme/theme1.xml, used by the one slide
# master / all layouts / all slides) is repainted with the default
# Office theme palette. Font scheme and format scheme (fills / lines /
# effects) are identical between the two themes, so only the 12 theme
# colour slots need to change.
#
# PowerPoint's object model exposes the live theme colour scheme via
# Slide.ThemeColorScheme (any slide works - they all share the single
# deck theme), with colours addressed in the standard MSO theme colour
# order: 1 Dark1, 2 Light1, 3 Dark2, 4 Light2, 5 Accent1 .. 10 Accent6,
# 11 Hyperlink, 12 FollowedHyperlink.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToRgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$tcs.Colors(1).RGB  = ToRgb 0x00 0x00 0x00   # Dark 1      -> 000000
$tcs.Colors(2).RGB  = ToRgb 0xFF 0xFF 0xFF   # Light 1     -> FFFFFF
$tcs.Colors(3).RGB  = ToRgb 0x44 0x54 0x6A   # Dark 2      -> 44546A
$tcs.Colors(4).RGB  = ToRgb 0xE7 0xE6 0xE6   # Light 2     -> E7E6E6
$tcs.Colors(5).RGB  = ToRgb 0x5B 0x9B 0xD5   # Accent 1    -> 5B9BD5
$tcs.Colors(6).RGB  = ToRgb 0xED 0x7D 0x31   # Accent 2    -> ED7D31
$tcs.Colors(7).RGB  = ToRgb 0xA5 0xA5 0xA5   # Accent 3    -> A5A5A5
$tcs.Colors(8).RGB  = ToRgb 0xFF 0xC0 0x00   # Accent 4    -> FFC000
$tcs.Colors(9).RGB  = ToRgb 0x44 0x72 0xC4   # Accent 5    -> 4472C4
$tcs.Colors(10).RGB = ToRgb 0x70 0xAD 0x47   # Accent 6    -> 70AD47
$tcs.Colors(11).RGB = ToRgb 0x05 0x63 0xC1   # Hyperlink   -> 0563C1
$tcs.Colors(12).RGB = ToRgb 0x95 0x4F 0x72   # Followed Hyperlink -> 954F72
